$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.81"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("D3").Value = "'41.08"
$ws.Range("E3").Value = "'0.06%"
$ws.Range("D4").Value = "'5.249"
$ws.Range("E4").Value = "'2.43%"
$ws.Range("D5").Value = "'0.07668"
$ws.Range("E5").Value = "'0.59%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.323"
$ws.Range("E6").Value = "'1.40%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.626"
$ws.Range("E7").Value = "'0.56%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9179"
$ws.Range("E8").Value = "'1.87%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.427"
$ws.Range("E9").Value = "'-3.21%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1226"
$ws.Range("E10").Value = "'12.14%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1831"
$ws.Range("E11").Value = "'3.31%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09107"
$ws.Range("E12").Value = "'-0.68%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04247"
$ws.Range("E13").Value = "'1.11%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1052"
$ws.Range("E14").Value = "'0.09%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001263"
$ws.Range("E15").Value = "'0.69%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005819"
$ws.Range("E16").Value = "'0.26%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007509"
$ws.Range("E17").Value = "'2,395.62%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.346"
$ws.Range("E18").Value = "'-0.34%"
$ws.Range("D20").Value = "'7.315"
$ws.Range("E20").Value = "'11.34%"
$ws.Range("E21").Value = "'1.79%"
$ws.Range("D23").Value = "'0.04065"
$ws.Range("E23").Value = "'0.00%"
$ws.Range("D24").Value = "'0.001264"
$ws.Range("D25").Value = "'0.004348"
$ws.Range("E25").Value = "'8.69%"
$ws.Range("E26").Value = "'-2.17%"
$ws.Range("D38").Value = "'0.02467"
$ws.Range("E38").Value = "'3.52%"
$ws.Range("D39").Value = "'0.05305"
$ws.Range("E39").Value = "'2.46%"
$ws.Range("D40").Value = "'0.007846"
$ws.Range("E40").Value = "'1.22%"
$ws.Range("D41").Value = "'0.1314"
$ws.Range("E41").Value = "'1.12%"
$ws.Range("D42").Value = "'0.006672"
$ws.Range("E42").Value = "'-2.43%"
$ws.Range("D43").Value = "'0.001914"
$ws.Range("E43").Value = "'-1.92%"
$ws.Range("D44").Value = "'0.007665"
$ws.Range("E44").Value = "'-10.31%"
$ws.Range("E45").Value = "'-0.54%"
$ws.Range("D46").Value = "'0.00006727"
$ws.Range("E46").Value = "'-3.75%"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("E48").Value = "'2,059.52%"
$ws.Range("E49").Value = "'-2.47%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.14%"
